$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VillaTest")
$ws.Activate()

# Update the hotel name in B2 from "Airport Hotel Earth" to "Hotel the View"
$ws.Range("B2").Value = "Hotel the View"

# Update the selection to B3 (matches resulting selection in the diff)
$ws.Range("B3").Select()
